$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" (last sheet) ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "PO Forecast"

# --- Header row ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$newSheet.Cells.Item(2,1).Value = 44948.99999999999
$newSheet.Cells.Item(2,2).Value = 11
$newSheet.Cells.Item(2,3).Value = -50.71503522655119
$newSheet.Cells.Item(2,4).Value = 72.20601510546081
$newSheet.Cells.Item(3,1).Value = 44983.99999999999
$newSheet.Cells.Item(3,2).Value = 20
$newSheet.Cells.Item(3,3).Value = -40.77122000448116
$newSheet.Cells.Item(3,4).Value = 78.69156784096532
$newSheet.Cells.Item(4,1).Value = 44997.99999999999
$newSheet.Cells.Item(4,2).Value = 24
$newSheet.Cells.Item(4,3).Value = -39.4520754078606
$newSheet.Cells.Item(4,4).Value = 89.21584179348331
$newSheet.Cells.Item(5,1).Value = 45004.99999999999
$newSheet.Cells.Item(5,2).Value = 26
$newSheet.Cells.Item(5,3).Value = -35.2771059326212
$newSheet.Cells.Item(5,4).Value = 85.65197899122418
$newSheet.Cells.Item(6,1).Value = 45032.99999999999
$newSheet.Cells.Item(6,2).Value = 33
$newSheet.Cells.Item(6,3).Value = -30.31262631849122
$newSheet.Cells.Item(6,4).Value = 99.3401553946444
$newSheet.Cells.Item(7,1).Value = 45060.99999999999
$newSheet.Cells.Item(7,2).Value = 41
$newSheet.Cells.Item(7,3).Value = -21.01152281342238
$newSheet.Cells.Item(7,4).Value = 108.4795274356375
$newSheet.Cells.Item(8,1).Value = 45067.99999999999
$newSheet.Cells.Item(8,2).Value = 42
$newSheet.Cells.Item(8,3).Value = -19.91312310361597
$newSheet.Cells.Item(8,4).Value = 110.4734565007998
$newSheet.Cells.Item(9,1).Value = 45074.99999999999
$newSheet.Cells.Item(9,2).Value = 44
$newSheet.Cells.Item(9,3).Value = -16.80105866065926
$newSheet.Cells.Item(9,4).Value = 110.2718335697239
$newSheet.Cells.Item(10,1).Value = 45081.99999999999
$newSheet.Cells.Item(10,2).Value = 46
$newSheet.Cells.Item(10,3).Value = -18.87989182602409
$newSheet.Cells.Item(10,4).Value = 107.3919922887264
$newSheet.Cells.Item(11,1).Value = 45095.99999999999
$newSheet.Cells.Item(11,2).Value = 50
$newSheet.Cells.Item(11,3).Value = -10.11865238430343
$newSheet.Cells.Item(11,4).Value = 109.1706227863063
$newSheet.Cells.Item(12,1).Value = 45102.99999999999
$newSheet.Cells.Item(12,2).Value = 52
$newSheet.Cells.Item(12,3).Value = -14.3895588939119
$newSheet.Cells.Item(12,4).Value = 118.5518673672159
$newSheet.Cells.Item(13,1).Value = 45109.99999999999
$newSheet.Cells.Item(13,2).Value = 54
$newSheet.Cells.Item(13,3).Value = -8.785054175619395
$newSheet.Cells.Item(13,4).Value = 116.2073025851737
$newSheet.Cells.Item(14,1).Value = 45137.99999999999
$newSheet.Cells.Item(14,2).Value = 61
$newSheet.Cells.Item(14,3).Value = 1.415408330803851
$newSheet.Cells.Item(14,4).Value = 124.8935133247708
$newSheet.Cells.Item(15,1).Value = 45144.99999999999
$newSheet.Cells.Item(15,2).Value = 63
$newSheet.Cells.Item(15,3).Value = 1.251061728979312
$newSheet.Cells.Item(15,4).Value = 124.8089314390538
$newSheet.Cells.Item(16,1).Value = 45158.99999999999
$newSheet.Cells.Item(16,2).Value = 67
$newSheet.Cells.Item(16,3).Value = 1.338613022166334
$newSheet.Cells.Item(16,4).Value = 129.6230439569366
$newSheet.Cells.Item(17,1).Value = 45165.99999999999
$newSheet.Cells.Item(17,2).Value = 68
$newSheet.Cells.Item(17,3).Value = 6.652228387373989
$newSheet.Cells.Item(17,4).Value = 129.022266056746
$newSheet.Cells.Item(18,1).Value = 45172.99999999999
$newSheet.Cells.Item(18,2).Value = 70
$newSheet.Cells.Item(18,3).Value = 6.922205646521072
$newSheet.Cells.Item(18,4).Value = 133.855641910626
$newSheet.Cells.Item(19,1).Value = 45179.99999999999
$newSheet.Cells.Item(19,2).Value = 72
$newSheet.Cells.Item(19,3).Value = 8.811962620163754
$newSheet.Cells.Item(19,4).Value = 134.3661575515019
$newSheet.Cells.Item(20,1).Value = 45186.99999999999
$newSheet.Cells.Item(20,2).Value = 74
$newSheet.Cells.Item(20,3).Value = 9.718937178510609
$newSheet.Cells.Item(20,4).Value = 135.2924551211508
$newSheet.Cells.Item(21,1).Value = 45193.99999999999
$newSheet.Cells.Item(21,2).Value = 76
$newSheet.Cells.Item(21,3).Value = 15.00699067726091
$newSheet.Cells.Item(21,4).Value = 136.7580072310625
$newSheet.Cells.Item(22,1).Value = 45200.99999999999
$newSheet.Cells.Item(22,2).Value = 78
$newSheet.Cells.Item(22,3).Value = 11.28377930107977
$newSheet.Cells.Item(22,4).Value = 141.320784579213
$newSheet.Cells.Item(23,1).Value = 45207.99999999999
$newSheet.Cells.Item(23,2).Value = 80
$newSheet.Cells.Item(23,3).Value = 11.02808234821897
$newSheet.Cells.Item(23,4).Value = 138.2257010968777
$newSheet.Cells.Item(24,1).Value = 45214.99999999999
$newSheet.Cells.Item(24,2).Value = 81
$newSheet.Cells.Item(24,3).Value = 19.53684112434618
$newSheet.Cells.Item(24,4).Value = 147.2252653489633
$newSheet.Cells.Item(25,1).Value = 45221.99999999999
$newSheet.Cells.Item(25,2).Value = 83
$newSheet.Cells.Item(25,3).Value = 19.63024605270061
$newSheet.Cells.Item(25,4).Value = 146.9076636326295
$newSheet.Cells.Item(26,1).Value = 45228.99999999999
$newSheet.Cells.Item(26,2).Value = 85
$newSheet.Cells.Item(26,3).Value = 22.25225120431753
$newSheet.Cells.Item(26,4).Value = 149.269131450153
$newSheet.Cells.Item(27,1).Value = 45235.99999999999
$newSheet.Cells.Item(27,2).Value = 87
$newSheet.Cells.Item(27,3).Value = 25.04745065190174
$newSheet.Cells.Item(27,4).Value = 150.247143570241
$newSheet.Cells.Item(28,1).Value = 45242.99999999999
$newSheet.Cells.Item(28,2).Value = 89
$newSheet.Cells.Item(28,3).Value = 29.53638671660407
$newSheet.Cells.Item(28,4).Value = 155.1613714204499
$newSheet.Cells.Item(29,1).Value = 45249.99999999999
$newSheet.Cells.Item(29,2).Value = 91
$newSheet.Cells.Item(29,3).Value = 30.16401819607203
$newSheet.Cells.Item(29,4).Value = 157.7875945732718
$newSheet.Cells.Item(30,1).Value = 45256.99999999999
$newSheet.Cells.Item(30,2).Value = 93
$newSheet.Cells.Item(30,3).Value = 24.53361236236935
$newSheet.Cells.Item(30,4).Value = 152.8072050616907

# --- Formatting: match style of existing sheets ---
# Header row style (bold, bordered, centered) copied from Weekly Quantity A1:B1
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Date column style copied from Weekly Quantity A2
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A30").PasteSpecial(-4122)

$excel.CutCopyMode = $false

Write-Host "PO Forecast sheet created."
